# Scheduled-runner refresh of market-price columns (H:N) across several
# "Leve Profits" sheets. Each row's currentAveragePrice* / LevePrice* /
# LeveProfit* values are re-pulled and rewritten; a handful of rows on
# the ALC sheet (125-141) no longer resolve to market data and have
# their H:N cells cleared entirely.

$wb = $excel.ActiveWorkbook

function Set-RowHN {
    param($ws, [int]$row, [object[]]$values)
    # $values must contain exactly 7 entries for columns H..N (use $null to clear a cell)
    $arr = New-Object 'object[,]' 1,7
    for ($i = 0; $i -lt 7; $i++) { $arr[0, $i] = $values[$i] }
    $ws.Range("H$row`:N$row").Value = $arr
}

# ---- ALC --------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 99: refreshed prices; LeveProfitNQ (M) no longer computable -> cleared
Set-RowHN $ws 99 @(4993.3335, 0, 4993.3335, 0, 14980.0005, $null, -17976.0005)

# Rows 125-141: no current market data -> clear the whole H:N block
$ws.Range("H125:N141").ClearContents()

# ---- ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
Set-RowHN $ws 2   @(2073.682, 2302.1, 1883.3334, 2302.1, 1883.3334, -2189.1, -2109.3334)
Set-RowHN $ws 116 @(2073.682, 2302.1, 1883.3334, 2302.1, 1883.3334, -8.099999999999909, -6471.3334)

# ---- BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
Set-RowHN $ws 3 @(2073.682, 2302.1, 1883.3334, 2302.1, 1883.3334, -2188.1, -2111.3334)

# ---- CRP ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
Set-RowHN $ws 16  @(777.58826, 701.63635, 916.8333, 701.63635, 916.8333, -414.63635, -1490.8333)
Set-RowHN $ws 31  @(43480424, 250001040, 2401.3157, 250001040, 2401.3157, -250000745, -2991.3157)
Set-RowHN $ws 34  @(43480424, 250001040, 2401.3157, 250001040, 2401.3157, -250000838, -2805.3157)
Set-RowHN $ws 113 @(777.58826, 701.63635, 916.8333, 701.63635, 916.8333, 1468.36365, -5256.8333)

# ---- CUL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
Set-RowHN $ws 5   @(754.9375, 462.18182, 1399, 1386.54546, 4197, -1274.54546, -4421)
Set-RowHN $ws 86  @(485.7143, 466.66666, 500, 1399.99998, 1500, -213.9999800000001, -3872)
Set-RowHN $ws 89  @(485.7143, 466.66666, 500, 4199.99994, 4500, 1728.00006, -16356)

# Row 131: only currentAveragePrice/HQ + LevePriceHQ + LeveProfitHQ moved
$ws.Range("H131").Value = 890.76
$ws.Range("J131").Value = 891.4433
$ws.Range("L131").Value = 2674.3299
$ws.Range("N131").Value = -12754.3299

Set-RowHN $ws 135 @(754.9375, 462.18182, 1399, 4159.63638, 12591, -1624.63638, -17661)

# ---- GSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
Set-RowHN $ws 132 @(2621.2903, 2298.348, 3549.75, 6895.044, 10649.25, -4365.044, -15709.25)
